# daily auto push: 2026-02-28 18:48 UTC
# Insert two new data rows (2026/02/28 土 and 2026/03/01 日) right before the
# existing "2026/12/29" block, which starts at row 879. All rows from 879
# downward shift down by two (to 881..922).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 879..920 down by inserting two blank rows at row 879.
$ws.Rows.Item(879).Resize(2).Insert()

# New row 879: 2026/02/28, 土, 22, 37
# (force column A to stay text, not get reinterpreted as a date serial,
# then drop back to the default "Normal" style so no explicit style index
# is left behind on the cell)
$ws.Cells.Item(879, 1).NumberFormat = "@"
$ws.Cells.Item(879, 1).Value = "2026/02/28"
$ws.Cells.Item(879, 1).Style = "Normal"
$ws.Cells.Item(879, 2).Value = "土"
$ws.Cells.Item(879, 3).Value = 22
$ws.Cells.Item(879, 4).Value = 37

# New row 880: 2026/03/01, 日, 1, 37
$ws.Cells.Item(880, 1).NumberFormat = "@"
$ws.Cells.Item(880, 1).Value = "2026/03/01"
$ws.Cells.Item(880, 1).Style = "Normal"
$ws.Cells.Item(880, 2).Value = "日"
$ws.Cells.Item(880, 3).Value = 1
$ws.Cells.Item(880, 4).Value = 37
